$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sentence "biophysicist double helix." was previously split into 3 tokens
# (row 2-4). It is now re-split into smaller word-level tokens, and a new
# sentence "3D model human genome, discrimination." is appended as
# individual word tokens (rows 5-9), replacing the old 3-token version of
# that sentence (which previously lived in rows 3-4).

# Row 2
$ws.Range("A2").Value = "biophysicist "
$ws.Range("B2").Value = 463
$ws.Range("C2").Value = 669
$ws.Range("D2").Value = 105
$ws.Range("E2").Value = 23

# Row 3
$ws.Range("A3").Value = "double "
$ws.Range("B3").Value = 564
$ws.Range("C3").Value = 694
$ws.Range("D3").Value = 63
$ws.Range("E3").Value = 23

# Row 4
$ws.Range("A4").Value = "helix."
$ws.Range("B4").Value = 627
$ws.Range("C4").Value = 694
$ws.Range("D4").Value = 45
$ws.Range("E4").Value = 23

# Row 5
$ws.Range("A5").Value = "3D "
$ws.Range("B5").Value = 168
$ws.Range("C5").Value = 851.8
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 23

# Row 6
$ws.Range("A6").Value = "model "
$ws.Range("B6").Value = 198
$ws.Range("C6").Value = 851.8
$ws.Range("D6").Value = 58
$ws.Range("E6").Value = 23

# Row 7
$ws.Range("A7").Value = "human "
$ws.Range("B7").Value = 350
$ws.Range("C7").Value = 1176
$ws.Range("D7").Value = 64
$ws.Range("E7").Value = 23

# Row 8
$ws.Range("A8").Value = "genome, "
$ws.Range("B8").Value = 414
$ws.Range("C8").Value = 1176
$ws.Range("D8").Value = 80
$ws.Range("E8").Value = 23

# Row 9
$ws.Range("A9").Value = "discrimination. "
$ws.Range("B9").Value = 243
$ws.Range("C9").Value = 1242.4
$ws.Range("D9").Value = 128
$ws.Range("E9").Value = 23
